$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -21.932
$ws.Range("E3").Value = 16.442

$ws.Range("A21").Value = -19.936

$ws.Range("A23").Value = -19.945
$ws.Range("E24").Value = 16.539

$ws.Range("A25").Value = -21.78900000000001

$ws.Range("C27").Value = -13.069

$ws.Range("C31").Value = -13.163

$ws.Range("C39").Value = -12.858

$ws.Range("C48").Value = -11.432

$ws.Range("C51").Value = -11.133

$ws.Range("C52").Value = -11.601

$ws.Range("C55").Value = -13.175

$ws.Range("C56").Value = -13.537

$ws.Range("A57").Value = -22.07999999999999
$ws.Range("C57").Value = -13.813
$ws.Range("E57").Value = 16.469

$ws.Range("A59").Value = -22.407

$ws.Range("E61").Value = 16.625

$ws.Range("A69").Value = -21.703

$ws.Range("E70").Value = 17.687

$ws.Range("C73").Value = -12.515

$ws.Range("A79").Value = -21.005

$ws.Range("A83").Value = -21.938

$ws.Range("E86").Value = 16.554

$ws.Range("C89").Value = -11.264

$ws.Range("C90").Value = -12.482

$ws.Range("A93").Value = -21.459

$ws.Range("E98").Value = 16.114

$ws.Range("E100").Value = 16.809

$ws.Range("E102").Value = 16.513
